$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B.
# This shifts the old RawActivations/PercActivations/totalActivation
# columns from B/C/D to C/D/E.
$ws.Columns.Item(2).Insert()

# New header cell B1 = "segments", formatted like the other header cells.
$ws.Cells.Item(1, 2).Value = "segments"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the segment-name labels from column A into the new column B
# (plain, unstyled, like they were before in B2:D20), and put the
# 0-based numeric index into column A (keeping column A's styling).
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $name = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 2).ClearFormats()
    $ws.Cells.Item($row, 1).Value = $i
}
